# "reverted PDF file and changed randomize func"
#
# The sheet used to hold a batch of 5 generated sender email/password pairs.
# The edit trims that back down to a single real sender row, turns the
# sender's e-mail address into a clickable mailto hyperlink, and updates a
# couple of UI/view bits (selection, window placement).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 2-5, leaving only row 1 (dimension becomes A1:B1).
$ws.Rows("2:5").Delete()

# Put the new sender credentials into the single remaining row.
$ws.Range("A1").Value = "emmanuelturner41@gmail.com"
$ws.Range("B1").Value = 'wwkkk#$1234'

# Turn the e-mail address into a clickable mailto: hyperlink. Excel
# automatically applies its built-in "Hyperlink" cell style (underlined,
# theme-colored font) to A1 when the link is added.
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:emmanuelturner41@gmail.com") | Out-Null

# Column A now holds a longer string than before, so let it resize to fit.
$ws.Columns("A:A").AutoFit() | Out-Null

# Match the saved selection / window placement recorded in the workbook.
$ws.Range("H8").Select() | Out-Null

$win = $excel.ActiveWindow
$win.Left = 1170
$win.Top = 1170
$win.Width = 21600
$win.Height = 11295
